# feat: inline entity hierarchy on event
#
# Adds a new "FinancialData -> address -> Address" property row to the
# Entities sheet (the FinancialData entity gains an inline `address`
# property of type Address), and leaves the workbook with the Entities
# sheet as the active/selected tab.

$wb = $excel.ActiveWorkbook

$entities = $wb.Worksheets.Item("Entities")
$rules = $wb.Worksheets.Item("Rules")

# New entity-property row describing FinancialData.address : Address
$entities.Range("B18").Value = "FinancialData"
$entities.Range("C18").Value = "address"
$entities.Range("D18").Value = "Address"

# Restore the Rules sheet's own last-known selection before switching away
# from it, then make Entities the active sheet/tab with its own selection.
$null = $rules.Range("C17").Select()

$null = $entities.Activate()
$null = $entities.Range("G24").Select()
